$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure D-column price cells keep their original text representation
# (the source data stores prices as text, e.g. "338.70", "28.212.77")
# by forcing text number-format before assigning number-like strings,
# so Excel does not silently convert them to floating point numbers
# and drop significant trailing zeros / reformat them.

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.212.77'
$ws.Range('E2').Value = '  +1.33%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.818.55'
$ws.Range('E3').Value = '  +2.56%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.75%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.70'
$ws.Range('E5').Value = '  -0.14%  '

# Row 6
$ws.Range('E6').Value = '  -0.56%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4391'
$ws.Range('E7').Value = '  +14.93%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3538'
$ws.Range('E8').Value = '  +3.92%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.70'
$ws.Range('E9').Value = '  -1.51%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.158'
$ws.Range('E10').Value = '  +1.75%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07461'
$ws.Range('E11').Value = '  +1.16%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.04'
$ws.Range('E12').Value = '  -1.69%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.9989'
$ws.Range('E13').Value = '  -0.59%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.285'
$ws.Range('E14').Value = '  -0.94%  '

# Row 15
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.309'
$ws.Range('E15').Value = '  -1.26%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.820.67'
$ws.Range('E16').Value = '  +2.64%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001088'
$ws.Range('E17').Value = '  +1.59%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06689'
$ws.Range('E18').Value = '  +0.41%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.21'

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9979'
$ws.Range('E20').Value = '  -0.72%  '

# Row 21
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.34'
$ws.Range('E21').Value = '  -0.45%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.469'
$ws.Range('E22').Value = '  +0.68%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.235.55'
$ws.Range('E23').Value = '  +1.29%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.15'
$ws.Range('E24').Value = '  +0.97%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.384'
$ws.Range('E25').Value = '  -0.70%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.492'
$ws.Range('E26').Value = '  +3.89%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.83'
$ws.Range('E27').Value = '  +0.89%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.54'
$ws.Range('E28').Value = '  +2.07%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.026.83'
$ws.Range('E29').Value = '  +2.98%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.310'
$ws.Range('E30').Value = '  -12.23%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '133.45'
$ws.Range('E31').Value = '  -0.21%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.068'
$ws.Range('E32').Value = '  +0.67%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.995'
$ws.Range('E33').Value = '  -0.37%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.09414'
$ws.Range('E34').Value = '  +6.13%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.38'
$ws.Range('E35').Value = '  -1.97%  '

# Row 36
$ws.Range('B36').Value = 'TheSandbox'
$ws.Range('C36').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6802'
$ws.Range('E36').Value = '  -0.04%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02373'
$ws.Range('E37').Value = '  -0.59%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.246'
$ws.Range('E38').Value = '  -0.53%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06275'
$ws.Range('E39').Value = '  -1.33%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2170'
$ws.Range('E40').Value = '  +0.55%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.486'
$ws.Range('E41').Value = '  -0.80%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.221'
$ws.Range('E42').Value = '  -0.29%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.300'
$ws.Range('E43').Value = '  +1.59%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9973'
$ws.Range('E44').Value = '  -0.64%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.06'
$ws.Range('E45').Value = '  -0.53%  '

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6175'
$ws.Range('E46').Value = '  -0.65%  '

# Row 47
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.878'
$ws.Range('E47').Value = '  +0.44%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '129.63'
$ws.Range('E48').Value = '  -2.62%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.050'
$ws.Range('E49').Value = '  -0.49%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.178'
$ws.Range('E50').Value = '  -2.27%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07121'
$ws.Range('E51').Value = '  -3.38%  '
